# Commit: Added test run for 5/7/2020
#
# The "Test Checklist" sheet (2nd worksheet) tracks one test run per column,
# starting at column C. The most recent run lives in column G (header date
# 4/30/2020, i.e. serial 43951) and every checklist row underneath it holds
# "PASS". This change adds a brand new run column (H) dated 5/7/2020 with
# the same "PASS" status in every row, and leaves the view scrolled/selected
# where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Checklist")
$ws.Activate()

$lastRow = 64

# Duplicate the most recent run column (G) into the new column (H) so the
# new column inherits the exact same cell styles (date format on the header,
# top-aligned text for every "PASS" cell below it) as well as the values.
$srcCol = $ws.Range("G1:G" + $lastRow)
$dstCol = $ws.Range("H1:H" + $lastRow)
$srcCol.Copy($dstCol)

# Re-point the header of the new column to this test run's date (5/7/2020).
# Using the date serial keeps the existing date-formatted style on H1 intact
# (assigning a date literal/string here would otherwise create a brand new
# number format instead of reusing style already used by C1:G1).
$ws.Range("H1").Value = 43958

# Restore the view to where the author left off: scrolled down so row 7 is
# at the top, with cell I22 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$ws.Range("I22").Select()
